# "New updates on Inheritance" - fill in the Topic Covered entries for
# the three rows that were still blank (dates 45352-45354 / C21:C23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = "Partial Classes"
$ws.Range("C22").Value = "Saturday : Holiday"
$ws.Range("C23").Value = "Static Classes, Static Properties, Inheritance"
